$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header strings: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 2) Turn the used range into an Excel table (ListObject) ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U55"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3) Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
